# Trainee Detail - New Team
# Update the "Group" sheet so each group's Project Code / Topic Code /
# Topic Name / Group Description columns are unique per group (1..4)
# instead of sharing the same generic text, and reorder + relabel the
# "Student" sheet's Group-Leader marker column.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Group": make project/topic/description text unique per row ----
$wsGroup = $wb.Worksheets.Item("Group")

$groupData = @{
    2 = @("Project Code 1", "Topic Code1", "Topic Name1", "Group Description1")
    3 = @("Project Code 2", "Topic Code2", "Topic Name2", "Group Description2")
    4 = @("Project Code 3", "Topic Code3", "Topic Name3", "Group Description3")
    5 = @("Project Code 4", "Topic Code4", "Topic Name4", "Group Description4")
}

foreach ($row in $groupData.Keys) {
    $vals = $groupData[$row]
    $wsGroup.Cells.Item($row, 3).Value = $vals[0]
    $wsGroup.Cells.Item($row, 4).Value = $vals[1]
    $wsGroup.Cells.Item($row, 5).Value = $vals[2]
    $wsGroup.Cells.Item($row, 6).Value = $vals[3]
}

# ---- Sheet "Student": reorder rows by group number and clear the leader flags ----
$wsStudent = $wb.Worksheets.Item("Student")

$studentData = @{
    2 = @(1, "nguyenvanc@gmail.com", "Nguyễn Văn C")
    3 = @(2, "nguyenvanh@gmail.com", "Nguyễn Văn H")
    4 = @(3, "nguyenvanl@gmail.com", "Nguyễn Văn L")
    5 = @(4, "anivns.com@gmail.com", "NVCK2002")
    6 = @(1, "nguyenvana@gmail.com", "Nguyen Van A")
    7 = @(2, "nguyenvanb@gmail.com", "Nguyen Van B")
    8 = @(3, "nguyenvank@gmail.com", "Nguyen Van K")
}

foreach ($row in $studentData.Keys) {
    $vals = $studentData[$row]
    $wsStudent.Cells.Item($row, 1).Value = $vals[0]
    $wsStudent.Cells.Item($row, 2).Value = $vals[1]
    $wsStudent.Cells.Item($row, 3).Value = $vals[2]
    $wsStudent.Cells.Item($row, 4).Value = ""
}

$wsStudent.Range("D10").Select() | Out-Null
